$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0154185022026432
$ws.Range("C2").Value = 0.959251101321586
$ws.Range("D2").Value = 0.0110132158590308
$ws.Range("E2").Value = 0.766519823788546
$ws.Range("F2").Value = 0.0253303964757709
$ws.Range("G2").Value = 0.0066079295154185
$ws.Range("H2").Value = 0.0220264317180617
$ws.Range("I2").Value = 0.861233480176211
$ws.Range("J2").Value = 0.00770925110132159
$ws.Range("K2").Value = 0.00330396475770925
$ws.Range("L2").Value = 0.0330396475770925
$ws.Range("M2").Value = 0.893171806167401
$ws.Range("N2").Value = 0.00110132158590308
$ws.Range("O2").Value = 0.00220264317180617
$ws.Range("P2").Value = 0.00330396475770925
$ws.Range("Q2").Value = 0.00550660792951542
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0.00220264317180617
$ws.Range("T2").Value = 0.00220264317180617
$ws.Range("U2").Value = 0.859030837004405
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0.0066079295154185
$ws.Range("X2").Value = 0.00550660792951542
$ws.Range("B3").Value = 0.00110132158590308
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.0275330396475771
$ws.Range("E3").Value = 0.0209251101321586
$ws.Range("F3").Value = 0.155286343612335
$ws.Range("G3").Value = 0.966960352422907
$ws.Range("H3").Value = 0.977973568281938
$ws.Range("I3").Value = 0.105726872246696
$ws.Range("J3").Value = 0.924008810572687
$ws.Range("K3").Value = 0.00991189427312775
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0.00110132158590308
$ws.Range("N3").Value = 0.865638766519824
$ws.Range("O3").Value = 0.965859030837004
$ws.Range("P3").Value = 0.0319383259911894
$ws.Range("Q3").Value = 0.00110132158590308
$ws.Range("R3").Value = 0.0187224669603524
$ws.Range("S3").Value = 0.00440528634361234
$ws.Range("T3").Value = 0.00110132158590308
$ws.Range("U3").Value = 0.0209251101321586
$ws.Range("V3").Value = 0.0066079295154185
$ws.Range("W3").Value = 0.0253303964757709
$ws.Range("X3").Value = 0.00220264317180617
$ws.Range("B4").Value = 0.975770925110132
$ws.Range("C4").Value = 0.0352422907488987
$ws.Range("D4").Value = 0.00110132158590308
$ws.Range("E4").Value = 0.0242290748898678
$ws.Range("F4").Value = 0.0308370044052863
$ws.Range("G4").Value = 0.0209251101321586
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.00440528634361234
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.00220264317180617
$ws.Range("L4").Value = 0.95704845814978
$ws.Range("M4").Value = 0.00330396475770925
$ws.Range("N4").Value = 0.122246696035242
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0.00550660792951542
$ws.Range("R4").Value = 0.0231277533039648
$ws.Range("S4").Value = 0.990088105726872
$ws.Range("T4").Value = 0.988986784140969
$ws.Range("U4").Value = 0.111233480176211
$ws.Range("V4").Value = 0.0242290748898678
$ws.Range("W4").Value = 0.964757709251101
$ws.Range("X4").Value = 0.964757709251101
$ws.Range("B5").Value = 0.00770925110132159
$ws.Range("C5").Value = 0.00550660792951542
$ws.Range("D5").Value = 0.960352422907489
$ws.Range("E5").Value = 0.188325991189427
$ws.Range("F5").Value = 0.788546255506608
$ws.Range("G5").Value = 0.00440528634361234
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.0286343612334802
$ws.Range("J5").Value = 0.0682819383259912
$ws.Range("K5").Value = 0.984581497797357
$ws.Range("L5").Value = 0.00991189427312775
$ws.Range("M5").Value = 0.102422907488987
$ws.Range("N5").Value = 0.0110132158590308
$ws.Range("O5").Value = 0.0319383259911894
$ws.Range("P5").Value = 0.964757709251101
$ws.Range("Q5").Value = 0.987885462555066
$ws.Range("R5").Value = 0.958149779735683
$ws.Range("S5").Value = 0.00330396475770925
$ws.Range("T5").Value = 0.00770925110132159
$ws.Range("U5").Value = 0.00881057268722467
$ws.Range("V5").Value = 0.969162995594714
$ws.Range("W5").Value = 0.00220264317180617
$ws.Range("X5").Value = 0.0253303964757709
